# Generate Report for handoff
# - "Handoff failed" -> "Not yet handed off" (Status column, row 2, all sheets)
# - Add "Latest Handoff File" entries (column C, row 2) for zh-cn / de-de, with hyperlinks
# - Update "Latest Handoff Datetime" (column D, row 2) with real timestamps
# - Update "Handoff Reason" (column H, row 2) from "Ignored" to "Include"

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/181437edecbab4a3c4b415f2693dce4df8cc06c1"
$srcId   = "aa6ca3cb-c3ce-4575-8446-36a6625c758e"
$fpHash  = "6ac4514a0b709804427a9b1ad8ec912cc1876bb2"

# ---- Status text: "Handoff failed" -> "Not yet handed off" everywhere it shows up ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = "Not yet handed off"
$wsOverview.Range("C2").Value2 = "Not yet handed off"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value2 = "Not yet handed off"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value2 = "Not yet handed off"

# ---- zh-cn sheet: populate the handoff-file / datetime / reason columns ----
$zhFile = "$srcId.$fpHash.zh-cn.xlf"
$zhUrl  = "$baseUrl/localization/zh-cn/$zhFile"

# Capture the existing A3 (.localization-config) hyperlink so it can be
# re-created *after* the new C2 one - this keeps the <hyperlinks> element
# in the same top-to-bottom, left-to-right order Excel would naturally use
# (A2, C2, A3) instead of appending C2 at the end of the list.
$zhA3Url = ""
$zhA3Disp = ""
$zhA3Link = $null
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$3") {
        $zhA3Url = $hl.Address
        $zhA3Disp = $hl.TextToDisplay
        $zhA3Link = $hl
    }
}
$zhA3Link.Delete()

$wsZh.Range("C2").Value2 = $zhFile
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhUrl, "", "", $zhFile) | Out-Null
$wsZh.Range("C2").Font.Underline = 2
$wsZh.Range("C2").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Url, "", "", $zhA3Disp) | Out-Null

$wsZh.Range("D2").Value2 = "2016-01-08 15:51:14"
$wsZh.Range("H2").Value2 = "Include"

# ---- de-de sheet: populate the handoff-file / datetime / reason columns ----
$deFile = "$srcId.$fpHash.de-de.xlf"
$deUrl  = "$baseUrl/localization/de-de/$deFile"

$deA3Url = ""
$deA3Disp = ""
$deA3Link = $null
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$3") {
        $deA3Url = $hl.Address
        $deA3Disp = $hl.TextToDisplay
        $deA3Link = $hl
    }
}
$deA3Link.Delete()

$wsDe.Range("C2").Value2 = $deFile
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deUrl, "", "", $deFile) | Out-Null
$wsDe.Range("C2").Font.Underline = 2
$wsDe.Range("C2").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Url, "", "", $deA3Disp) | Out-Null

$wsDe.Range("D2").Value2 = "2016-01-08 15:51:30"
$wsDe.Range("H2").Value2 = "Include"
